$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Insert three new blank columns before column B (the existing
#    B,C,D,E data - including the date headers in row 1 - shifts
#    right to E,F,G,H). This matches the report's "roll in a new
#    week of columns" pattern used for this MarketBeat rank sheet.
# ------------------------------------------------------------------
$ws.Range("B1:D1").EntireColumn.Insert(-4161)  # xlShiftToRight

# Keep the same narrow "8.0" column width used by the rest of the report
# on the three freshly-inserted columns (C, D, E) as well as the column
# that the old "E" data was pushed into (H).
$ws.Columns.Item(3).ColumnWidth = 7.1666666666667
$ws.Columns.Item(4).ColumnWidth = 7.1666666666667
$ws.Columns.Item(5).ColumnWidth = 7.1666666666667
$ws.Columns.Item(6).ColumnWidth = 7.1666666666667
$ws.Columns.Item(7).ColumnWidth = 7.1666666666667
$ws.Columns.Item(8).ColumnWidth = 7.1666666666667

# ------------------------------------------------------------------
# 2. New header dates for the freshly inserted columns.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# ------------------------------------------------------------------
# 3. Fill the new B,C,D columns with the "UN" placeholder for every
#    data row, matching the rest of the grid.
# ------------------------------------------------------------------
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("B$r").Value = "UN"
    $ws.Range("C$r").Value = "UN"
    $ws.Range("D$r").Value = "UN"
}

# ------------------------------------------------------------------
# 4. Overwrite specific cells with real MarketBeat rating updates.
# ------------------------------------------------------------------
$ws.Range("C18").Value = "6/18/2018,Reiterates,Buy,`$87.00"
$ws.Range("D18").Value = "6/18/2018,Reiterates,Buy,`$87.00"

$ws.Range("B20").Value = "6/26/2018,Reiterates,Buy -> Buy,`$88.00 -> `$98.00"
$ws.Range("C20").Value = "6/26/2018,Reiterates,Buy -> Buy,`$88.00 -> `$98.00"
$ws.Range("D20").Value = "6/26/2018,Reiterates,Buy -> Buy,`$88.00 -> `$98.00"

# ------------------------------------------------------------------
# 5. Append a new analyst-coverage group (two new rows).
# ------------------------------------------------------------------
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
